# The workbook tracks two related "transition matrix" sheets: a header
# row (B1:G1) and row labels (A2:A30) built from state-name combinations.
# The state that used to be called "ScreenRecStarted" was renamed to
# "0_unstated" - update the header cell and the four row labels that are
# built from it. The other (unchanged) labels keep their original text;
# Excel will naturally collapse/re-point the shared-string table so the
# unused "ScreenRecStarted" string disappears on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: the last column used to read "ScreenRecStarted".
$ws.Range("G1").Value = "0_unstated"

# Row labels that combined the renamed state with the other states.
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Move/update the saved selection to match the author's final cursor
# position.
$ws.Range("E14").Select()
